$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column D ("enlace_de_drive") entirely, shifting E:H left to D:G
$ws.Range("D1:D3").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

# Update row 2 with the new record's data (columns now: A,B,C,D,E,F,G)
$ws.Range("A2").Value = "dmarmols@miumg.edu.gt"

# B2 must stay a text string (preserve the leading zero), so force text format
# before assigning, then restore the default (unstyled) cell style.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "09072313365"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "43aa90963efd17f41ea64c8a2e8ded98"
$ws.Range("F2").Value = 4
